# --- Applies the "Update 02. - 05. 9." commit: appends 23 new documentation
# --- rows (2018-09-02 .. 2018-09-05 drone build/flight log entries).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Dates (column A) for the new rows.
$ws.Range("A240").Value = 43345
$ws.Range("A241").Value = 43345
$ws.Range("A242").Value = 43345
$ws.Range("A243").Value = 43345
$ws.Range("A244").Value = 43345
$ws.Range("A245").Value = 43345
$ws.Range("A246").Value = 43347
$ws.Range("A247").Value = 43347
$ws.Range("A248").Value = 43347
$ws.Range("A249").Value = 43347
$ws.Range("A250").Value = 43347
$ws.Range("A251").Value = 43347
$ws.Range("A252").Value = 43347
$ws.Range("A253").Value = 43347
$ws.Range("A254").Value = 43348
$ws.Range("A255").Value = 43348
$ws.Range("A256").Value = 43348
$ws.Range("A257").Value = 43348
$ws.Range("A258").Value = 43348
$ws.Range("A259").Value = 43348
$ws.Range("A260").Value = 43348
$ws.Range("A261").Value = 43348
$ws.Range("A262").Value = 43348

# Filenames / descriptions (columns B, C), written in original entry order.
$ws.Range("B240").Value = '2018-09-04 1.JPG'
$ws.Range("C240").Value = 'Nach einem Flugversuch war der kontroller vollkommen entstellt. Nach dem Erfolg ist die konstruktion zu instabil. Da muss was neues her'
$ws.Range("B241").Value = '2018-09-02 6.JPG'
$ws.Range("C241").Value = 'Zunächst muss die powerbank weg. Jetzt wurde eine andere aufgebrochen und die Batterie unter der Drohne befestigt'
$ws.Range("B242").Value = '2018-09-02 7.JPG'
$ws.Range("C242").Value = 'Das Kontroll-Brett mit Ladebuchse, Stromversorgung, Schalter und Ladungsstatus kommt als neues BMS an die alte stelle'
$ws.Range("B243").Value = '2018-09-02 8.JPG'
$ws.Range("C243").Value = 'Hier der Anblick von hinten. Die drohne ist jetzt viel flacher, und die Sensoren sitzen näher am Masse/ Steuerzentrum. Zum an und ausschalten gibt es noch einen Kippschalter'
$ws.Range("B244").Value = '2018-09-02 9.JPG'
$ws.Range("C244").Value = 'Von vorne… Die Drohne ist schon flacher, aber die Schaltbretter müssen noch gerade befestigt werden'
$ws.Range("B245").Value = '2018-09-02 10.JPG'
$ws.Range("C245").Value = 'Die lösung auf das problem mit den vielen Steckern vom Rahmen zum Controller wurde jetzt in die Hand genommen und zusammengefasst'
$ws.Range("B246").Value = '2018-09-04 1.JPG'
$ws.Range("C246").Value = 'Die Batterie hat jetzt auch ein sicheres Abteil und die Kamera einen Platz…'
$ws.Range("B247").Value = '2018-09-04 2.JPG'
$ws.Range("C247").Value = 'Das Gesicht der neuen Drohne'
$ws.Range("B248").Value = '2018-09-04 3.JPG'
$ws.Range("C248").Value = 'Die Sicht der Drohne'
$ws.Range("B249").Value = '2018-09-04 4.JPG'
$ws.Range("C249").Value = 'Trotz feuchtem Wetter muss mit dem neuen aufbau ein neuer Testflug geflogen werden'
$ws.Range("B251").Value = '2018-09-04 6.JPG'
$ws.Range("C250").Value = 'Patrick ist da, um schöne Bilder zu machen. Diesmal wird es sogar Hochaufgelöst'
$ws.Range("B250").Value = '2018-09-04 5.JPG'
$ws.Range("C251").Value = 'Aufschalten auf die Drohne'
$ws.Range("B252").Value = '2018-09-04 7.MOV'
$ws.Range("C252").Value = 'Die Drohne hebt nicht ab. Es waren 2 Propeller falsch rum draufgesessen'
$ws.Range("B253").Value = '2018-09-04 8.MOV'
$ws.Range("C253").Value = 'Schon  wieder ein absturz. Das ist niederschmetternd. Nach dem erfolg vor 2 Tagen, dachte ich, wir könnten es schaffen. Jetzt bin ich Zeitmäßig wieder 2-3 Tage los.'
$ws.Range("C254").Value = 'Heute ist schönes Wetter. Zeit um dem Problemen von Gestern nachzugehen. Der Start ist gut, bloß dann wieder, ein Salto'
$ws.Range("C255").Value = 'Und wieder'
$ws.Range("C256").Value = 'Und wieder'
$ws.Range("C258").Value = 'Und wieder'
$ws.Range("B256").Value = '2018-09-05 3.mp4'
$ws.Range("B255").Value = '2018-09-05 2.AVI'
$ws.Range("B254").Value = '2018-09-05 1.AVI'
$ws.Range("B257").Value = '2018-09-05 4.mp4'
$ws.Range("C257").Value = 'Falsche Propeller'
$ws.Range("B258").Value = '2018-09-05 5.mp4'
$ws.Range("B259").Value = '2018-09-05 6.AVI'
$ws.Range("C259").Value = 'Dieser Flug war schon wieder Okay. So gut, wie man mit einer Zeitverzögerung von 1-2s halt fliegen kann. Fernsteuern ist langsam aussichtslos'
$ws.Range("B260").Value = '2018-09-05 7.AVI'
$ws.Range("C260").Value = 'Die Flüge werden immer besser. Ich kann die richtung einigermaßen bestimmen und die Drohne steht gerade in der Luft. Nur das reagieren auf das steigen und sinken der Drohne funktioniert nicht'
$ws.Range("B261").Value = '2018-09-05 8.mp4'
$ws.Range("C261").Value = 'Apropos steigen und sinken. Den Kirschbaum brauchen wir auch nicht mehr schneiden :)'
$ws.Range("B262").Value = '2018-09-05 9.mp4'
$ws.Range("C262").Value = 'Kontrollieren der richtungen in trocken+bungen. Soweit scheint alles zu klappen. Wüde die Drohne nicht abdriften, könnte man schon die erste automatation reinbringen'

# Copy the date-format (numFmt 14, m/d/yyyy) from the last pre-existing
# date cell (A239) onto the new date cells so they render as dates like
# the rest of column A.
$ws.Range("A239").Copy() | Out-Null
$ws.Range("A240:A262").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Leave the selection where the author left off editing.
$ws.Range("B263").Select() | Out-Null
